$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Extend the "Description" paragraph text.
# ---------------------------------------------------------------------------
$oldDescription = "The Arcade game based on movie Stuart Little (1999) in style of 70-80s Atari 2600 games. "
$newDescription = "The Arcade game based on movie Stuart Little (1999) in style of 70-80s Atari 2600 games. It has 3 levels, which are based on “survival highscore” mode. Still, some of them could be completed."

$findRange = $d.Content
$found = $findRange.Find.Execute($oldDescription, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $newDescription, 2)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark so it sits right after the text we just
#    inserted (end of the Description paragraph), instead of its original
#    spot in front of the "Assemgul' Amantaeva" paragraph.
#
#    A directly-collapsed Range placed exactly at a paragraph's last
#    character (just before the paragraph mark) cannot reliably anchor a
#    bookmark, so we temporarily insert a marker character right after the
#    target spot, anchor the bookmark there (now a genuine "inside text"
#    position), then remove the marker again - the bookmark stays put.
# ---------------------------------------------------------------------------
$endPos = $findRange.End

$marker = $d.Range($endPos, $endPos)
$marker.InsertAfter("X")

$bookmarkSpot = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$markerChar = $d.Range($endPos, $endPos + 1)
$markerChar.Delete()

# ---------------------------------------------------------------------------
# 3) Update Token Sabit's role from Tester to Programmer assistant.
# ---------------------------------------------------------------------------
$roleRange = $d.Content
$roleRange.Find.Execute("**Token Sabit - Tester", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**Token Sabit - Programmer assistant", 2) | Out-Null
